$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Update the "CasesTab" query text (row 2, column B) to drop the trailing
# "Cohort" column from the RETURN clause (coalesce(co.cohort_description, '') AS `Cohort`)
$newCasesQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`nWHERE s.clinical_study_designation IN ['GLIOMA01'] `nMATCH (c)<--(diag:diagnosis)`nOPTIONAL MATCH (samp:sample)-->(c)`nOPTIONAL MATCH (co:cohort)<-[*]-(c)`nWITH DISTINCT c, s, demo, diag, co`nRETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n        coalesce(demo.breed, '') AS Breed ,`n        coalesce(diag.disease_term, '') AS Diagnosis ,`n        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n        coalesce(demo.sex, '') AS Sex ,`n        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n        coalesce(demo.weight, '') AS ``Weight (kg)``,`n        coalesce(diag.best_response, '') AS ``Response to Treatment``"

$ws.Range("B2").Value = $newCasesQuery

# Selection moves to B2 and the view scrolls back to show row 1 (topLeftCell reset)
$ws.Range("B2").Select()
